$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Version value (row 3)
$ws.Range("B3").Value = "0.2.0"

# Update the Date value (row 8)
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# Insert a new row after "Contact" (row 10), before "Description" (row 11),
# shifting all subsequent rows down by one.
$ws.Rows.Item(11).Insert()

# Copy formatting (borders / fill / alignment) from the row above so the
# new row matches the rest of the table.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new "Jurisdiction" row.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
